$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 29999
$ws.Range("J13").Value = 29999
$ws.Range("L13").Value = 29999
$ws.Range("N13").Value = -30337
$ws.Range("H17").Value = 401870.7
$ws.Range("J17").Value = 441895.06
$ws.Range("L17").Value = 1325685.18
$ws.Range("N17").Value = -1326021.18
$ws.Range("H64").Value = 35717972
$ws.Range("I64").Value = 62502576
$ws.Range("J64").Value = 5166.6665
$ws.Range("K64").Value = 62502576
$ws.Range("L64").Value = 5166.6665
$ws.Range("M64").Value = -62502328
$ws.Range("N64").Value = -5662.6665
$ws.Range("H67").Value = 35717972
$ws.Range("I67").Value = 62502576
$ws.Range("J67").Value = 5166.6665
$ws.Range("K67").Value = 62502576
$ws.Range("L67").Value = 5166.6665
$ws.Range("M67").Value = -62501718
$ws.Range("N67").Value = -6882.6665
$ws.Range("H70").Value = 6671.357
$ws.Range("I70").Value = 2333.3333
$ws.Range("J70").Value = 7854.4546
$ws.Range("K70").Value = 6999.999899999999
$ws.Range("L70").Value = 23563.3638
$ws.Range("M70").Value = -6729.999899999999
$ws.Range("N70").Value = -24103.3638
$ws.Range("H73").Value = 6671.357
$ws.Range("I73").Value = 2333.3333
$ws.Range("J73").Value = 7854.4546
$ws.Range("K73").Value = 6999.999899999999
$ws.Range("L73").Value = 23563.3638
$ws.Range("M73").Value = -6063.999899999999
$ws.Range("N73").Value = -25435.3638
$ws.Range("H74").Value = 63628324
$ws.Range("I74").Value = 101801576
$ws.Range("J74").Value = 6233.3335
$ws.Range("K74").Value = 101801576
$ws.Range("L74").Value = 6233.3335
$ws.Range("M74").Value = -101800640
$ws.Range("N74").Value = -8105.3335
$ws.Range("H77").Value = 63628324
$ws.Range("I77").Value = 101801576
$ws.Range("J77").Value = 6233.3335
$ws.Range("K77").Value = 509007880
$ws.Range("L77").Value = 31166.6675
$ws.Range("M77").Value = -509003200
$ws.Range("N77").Value = -40526.6675
$ws.Range("H80").Value = 1141.4286
$ws.Range("I80").Value = 1493.2858
$ws.Range("K80").Value = 4479.857400000001
$ws.Range("M80").Value = -3481.857400000001
$ws.Range("H83").Value = 1141.4286
$ws.Range("I83").Value = 1493.2858
$ws.Range("K83").Value = 13439.5722
$ws.Range("M83").Value = -8447.572200000001
$ws.Range("H98").Value = 1156.6786
$ws.Range("I98").Value = 1156.6786
$ws.Range("K98").Value = 1156.6786
$ws.Range("M98").Value = 341.3214
$ws.Range("H106").Value = 3678.5
$ws.Range("I106").Value = 3912.8333
$ws.Range("J106").Value = 2975.5
$ws.Range("K106").Value = 3912.8333
$ws.Range("L106").Value = 2975.5
$ws.Range("M106").Value = -3281.8333
$ws.Range("N106").Value = -4237.5
$ws.Range("H107").Value = 878.6667
$ws.Range("I107").Value = 914.6
$ws.Range("J107").Value = 699
$ws.Range("K107").Value = 914.6
$ws.Range("L107").Value = 699
$ws.Range("M107").Value = 1005.4
$ws.Range("N107").Value = -4539
$ws.Range("H109").Value = 41666.332
$ws.Range("I109").Value = 29999
$ws.Range("K109").Value = 29999
$ws.Range("M109").Value = -28612
$ws.Range("H122").Value = 1156.6786
$ws.Range("I122").Value = 1156.6786
$ws.Range("K122").Value = 3470.0358
$ws.Range("M122").Value = -1020.0358
$ws.Range("H129").Value = 9906.1
$ws.Range("I129").Value = 1451.5555
$ws.Range("K129").Value = 4354.666499999999
$ws.Range("M129").Value = 645.3335000000006
$ws.Range("H137").Value = 1831.079
$ws.Range("I137").Value = 1440.6666
$ws.Range("J137").Value = 2789.3635
$ws.Range("K137").Value = 4321.9998
$ws.Range("L137").Value = 8368.0905
$ws.Range("M137").Value = -1771.9998
$ws.Range("N137").Value = -13468.0905

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 702098.1
$ws.Range("I2").Value = 920378.1
$ws.Range("J2").Value = 3602.2
$ws.Range("K2").Value = 920378.1
$ws.Range("L2").Value = 3602.2
$ws.Range("M2").Value = -920265.1
$ws.Range("N2").Value = -3828.2
$ws.Range("H32").Value = 1599.0212
$ws.Range("I32").Value = 1599.0212
$ws.Range("K32").Value = 1599.0212
$ws.Range("M32").Value = -1312.0212
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").Value = ""
$ws.Range("H45").Value = 2963.3333
$ws.Range("I45").Value = 1445
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 1445
$ws.Range("L45").Value = 6000
$ws.Range("M45").Value = -1068
$ws.Range("N45").Value = -6754
$ws.Range("H61").Value = 111114010
$ws.Range("I61").Value = 142860160
$ws.Range("J61").Value = 2498.5
$ws.Range("K61").Value = 142860160
$ws.Range("L61").Value = 2498.5
$ws.Range("M61").Value = -142859948
$ws.Range("N61").Value = -2922.5
$ws.Range("H74").Value = 25644474
$ws.Range("I74").Value = 31253026
$ws.Range("K74").Value = 31253026
$ws.Range("M74").Value = -31252152
$ws.Range("H77").Value = 25644474
$ws.Range("I77").Value = 31253026
$ws.Range("K77").Value = 156265130
$ws.Range("M77").Value = -156260762
$ws.Range("H110").Value = 335328
$ws.Range("I110").Value = 335328
$ws.Range("K110").Value = 335328
$ws.Range("M110").Value = -333283
$ws.Range("H111").Value = 58993
$ws.Range("J111").Value = 58993
$ws.Range("L111").Value = 58993
$ws.Range("N111").Value = -67173
$ws.Range("H116").Value = 702098.1
$ws.Range("I116").Value = 920378.1
$ws.Range("J116").Value = 3602.2
$ws.Range("K116").Value = 920378.1
$ws.Range("L116").Value = 3602.2
$ws.Range("M116").Value = -918084.1
$ws.Range("N116").Value = -8190.2
$ws.Range("H122").Value = 5744.619
$ws.Range("I122").Value = 4781.85
$ws.Range("K122").Value = 14345.55
$ws.Range("M122").Value = -11895.55
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = ""
$ws.Range("H132").Value = 2780955.8
$ws.Range("I132").Value = 2780955.8
$ws.Range("K132").Value = 8342867.399999999
$ws.Range("M132").Value = -8340337.399999999
$ws.Range("H136").Value = 111114010
$ws.Range("I136").Value = 142860160
$ws.Range("J136").Value = 2498.5
$ws.Range("K136").Value = 428580480
$ws.Range("L136").Value = 7495.5
$ws.Range("M136").Value = -428577930
$ws.Range("N136").Value = -12595.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 702098.1
$ws.Range("I3").Value = 920378.1
$ws.Range("J3").Value = 3602.2
$ws.Range("K3").Value = 920378.1
$ws.Range("L3").Value = 3602.2
$ws.Range("M3").Value = -920264.1
$ws.Range("N3").Value = -3830.2
$ws.Range("H20").Value = 1638.8334
$ws.Range("I20").Value = 2264.5
$ws.Range("J20").Value = 1513.7
$ws.Range("K20").Value = 2264.5
$ws.Range("L20").Value = 1513.7
$ws.Range("M20").Value = -2017.5
$ws.Range("N20").Value = -2007.7
$ws.Range("H86").Value = 5250.5
$ws.Range("I86").Value = 5250.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5250.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -4127.5
$ws.Range("N86").Value = ""
$ws.Range("H89").Value = 5250.5
$ws.Range("I89").Value = 5250.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 26252.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -20636.5
$ws.Range("N89").Value = ""
$ws.Range("H94").Value = 4811.7144
$ws.Range("I94").Value = 6521.923
$ws.Range("K94").Value = 6521.923
$ws.Range("M94").Value = -6070.923
$ws.Range("H99").Value = 2360.4783
$ws.Range("I99").Value = 2619.3
$ws.Range("J99").Value = 2161.3845
$ws.Range("K99").Value = 2619.3
$ws.Range("L99").Value = 2161.3845
$ws.Range("M99").Value = -1121.3
$ws.Range("N99").Value = -5157.3845
$ws.Range("H107").Value = 337399.66
$ws.Range("J107").Value = 504999.5
$ws.Range("L107").Value = 504999.5
$ws.Range("N107").Value = -508839.5
$ws.Range("H125").Value = 66890
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 66890
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 66890
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = -76730
$ws.Range("H134").Value = 12630445
$ws.Range("I134").Value = 13295024
$ws.Range("K134").Value = 39885072
$ws.Range("M134").Value = -39882537
$ws.Range("H137").Value = 50000
$ws.Range("J137").Value = 50000
$ws.Range("L137").Value = 50000
$ws.Range("N137").Value = -60200
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1037.4546
$ws.Range("I19").Value = 789
$ws.Range("K19").Value = 789
$ws.Range("M19").Value = -619
$ws.Range("H24").Value = 1037.4546
$ws.Range("I24").Value = 789
$ws.Range("K24").Value = 789
$ws.Range("M24").Value = -619
$ws.Range("H31").Value = 4507.724
$ws.Range("I31").Value = 3017.1765
$ws.Range("K31").Value = 3017.1765
$ws.Range("M31").Value = -2722.1765
$ws.Range("H32").Value = 21958.334
$ws.Range("I32").Value = 5875
$ws.Range("K32").Value = 5875
$ws.Range("M32").Value = -5559
$ws.Range("H34").Value = 4507.724
$ws.Range("I34").Value = 3017.1765
$ws.Range("K34").Value = 3017.1765
$ws.Range("M34").Value = -2815.1765
$ws.Range("H35").Value = 1905.7142
$ws.Range("I35").Value = 1946.6666
$ws.Range("K35").Value = 1946.6666
$ws.Range("M35").Value = -1652.6666
$ws.Range("H58").Value = 21751008
$ws.Range("I58").Value = 21751008
$ws.Range("K58").Value = 21751008
$ws.Range("M58").Value = -21750805
$ws.Range("H59").Value = 93332.664
$ws.Range("I59").Value = 10000
$ws.Range("J59").Value = 109999.2
$ws.Range("K59").Value = 10000
$ws.Range("L59").Value = 109999.2
$ws.Range("M59").Value = -8855
$ws.Range("N59").Value = -112289.2
$ws.Range("H99").Value = 3218.7778
$ws.Range("J99").Value = 3356.1428
$ws.Range("L99").Value = 3356.1428
$ws.Range("N99").Value = -6352.1428
$ws.Range("H105").Value = 5103766
$ws.Range("I105").Value = 5103766
$ws.Range("K105").Value = 5103766
$ws.Range("M105").Value = -5102019
$ws.Range("H107").Value = 758037.75
$ws.Range("I107").Value = 926642.9399999999
$ws.Range("K107").Value = 926642.9399999999
$ws.Range("M107").Value = -924722.9399999999
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H121").Value = 100000
$ws.Range("J121").Value = 100000
$ws.Range("L121").Value = 100000
$ws.Range("N121").Value = -102620
$ws.Range("H122").Value = 4489
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 4489
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 13467
$ws.Range("M122").Value = ""
$ws.Range("N122").Value = -18367
$ws.Range("H126").Value = 3218.7778
$ws.Range("J126").Value = 3356.1428
$ws.Range("L126").Value = 10068.4284
$ws.Range("N126").Value = -15008.4284
$ws.Range("H132").Value = 142860290
$ws.Range("I132").Value = 166669650
$ws.Range("J132").Value = 4100
$ws.Range("K132").Value = 500008950
$ws.Range("L132").Value = 12300
$ws.Range("M132").Value = -500006420
$ws.Range("N132").Value = -17360
$ws.Range("H136").Value = 21751008
$ws.Range("I136").Value = 21751008
$ws.Range("K136").Value = 65253024
$ws.Range("M136").Value = -65250474
$ws.Range("H137").Value = 172492.62
$ws.Range("I137").Value = 30000
$ws.Range("J137").Value = 192848.72
$ws.Range("K137").Value = 30000
$ws.Range("L137").Value = 192848.72
$ws.Range("M137").Value = -24900
$ws.Range("N137").Value = -203048.72
$ws.Range("H138").Value = 107500
$ws.Range("J138").Value = 107500
$ws.Range("L138").Value = 107500
$ws.Range("N138").Value = -117780
$ws.Range("H140").Value = 67246.25
$ws.Range("J140").Value = 67246.25
$ws.Range("L140").Value = 67246.25
$ws.Range("N140").Value = -77606.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 112198.89
$ws.Range("I5").Value = 200451.2
$ws.Range("K5").Value = 601353.6000000001
$ws.Range("M5").Value = -601241.6000000001
$ws.Range("H7").Value = 2500311
$ws.Range("I7").Value = 3333583
$ws.Range("J7").Value = 495
$ws.Range("K7").Value = 10000749
$ws.Range("L7").Value = 1485
$ws.Range("M7").Value = -10000637
$ws.Range("N7").Value = -1709
$ws.Range("H34").Value = 875.3333
$ws.Range("I34").Value = 896.5
$ws.Range("J34").Value = 833
$ws.Range("K34").Value = 2689.5
$ws.Range("L34").Value = 2499
$ws.Range("M34").Value = -2605.5
$ws.Range("N34").Value = -2667
$ws.Range("H39").Value = 2000.1428
$ws.Range("J39").Value = 4000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588
$ws.Range("H46").Value = 666.3333
$ws.Range("I46").Value = 749.5
$ws.Range("J46").Value = 500
$ws.Range("K46").Value = 2248.5
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -2157.5
$ws.Range("N46").Value = -1682
$ws.Range("H55").Value = 1250
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = ""
$ws.Range("H62").Value = 11506.5
$ws.Range("J62").Value = 11506.5
$ws.Range("L62").Value = 34519.5
$ws.Range("N62").Value = -35891.5
$ws.Range("H65").Value = 11506.5
$ws.Range("J65").Value = 11506.5
$ws.Range("L65").Value = 103558.5
$ws.Range("N65").Value = -110422.5
$ws.Range("H81").Value = 308062.25
$ws.Range("I81").Value = 308062.25
$ws.Range("K81").Value = 924186.75
$ws.Range("M81").Value = -923063.75
$ws.Range("H84").Value = 308062.25
$ws.Range("I84").Value = 308062.25
$ws.Range("K84").Value = 2772560.25
$ws.Range("M84").Value = -2766944.25
$ws.Range("H92").Value = 298
$ws.Range("I92").Value = 174.5
$ws.Range("J92").Value = 359.75
$ws.Range("K92").Value = 523.5
$ws.Range("L92").Value = 1079.25
$ws.Range("M92").Value = 724.5
$ws.Range("N92").Value = -3575.25
$ws.Range("H94").Value = 20891.777
$ws.Range("I94").Value = 9000
$ws.Range("J94").Value = 24289.428
$ws.Range("K94").Value = 27000
$ws.Range("L94").Value = 72868.284
$ws.Range("M94").Value = -26324
$ws.Range("N94").Value = -74220.284
$ws.Range("H98").Value = 498.25
$ws.Range("I98").Value = 447.66666
$ws.Range("J98").Value = 650
$ws.Range("K98").Value = 1342.99998
$ws.Range("L98").Value = 1950
$ws.Range("M98").Value = 155.0000199999999
$ws.Range("N98").Value = -4946
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 3000
$ws.Range("M110").Value = 1090
$ws.Range("H113").Value = 143508.72
$ws.Range("I113").Value = 333554.34
$ws.Range("J113").Value = 974.5
$ws.Range("K113").Value = 1000663.02
$ws.Range("L113").Value = 2923.5
$ws.Range("M113").Value = -998493.02
$ws.Range("N113").Value = -7263.5
$ws.Range("H135").Value = 112198.89
$ws.Range("I135").Value = 200451.2
$ws.Range("K135").Value = 1804060.8
$ws.Range("M135").Value = -1801525.8
$ws.Range("H139").Value = 877.25
$ws.Range("I139").Value = 877.25
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 2631.75
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = 2508.25
$ws.Range("N139").Value = ""
$ws.Range("H140").Value = 2065.7693
$ws.Range("I140").Value = 2065.7693
$ws.Range("K140").Value = 6197.3079
$ws.Range("M140").Value = -1017.3079
$ws.Range("H141").Value = 5594.5
$ws.Range("I141").Value = 1189
$ws.Range("K141").Value = 3567
$ws.Range("M141").Value = 1613

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 17450
$ws.Range("I40").Value = 17450
$ws.Range("K40").Value = 17450
$ws.Range("M40").Value = -17299
$ws.Range("H70").Value = 6510.778
$ws.Range("I70").Value = 7132.6665
$ws.Range("K70").Value = 7132.6665
$ws.Range("M70").Value = -6862.6665
$ws.Range("H73").Value = 6510.778
$ws.Range("I73").Value = 7132.6665
$ws.Range("K73").Value = 7132.6665
$ws.Range("M73").Value = -6196.6665
$ws.Range("H97").Value = 1549.4546
$ws.Range("I97").Value = 1572.2354
$ws.Range("K97").Value = 1572.2354
$ws.Range("M97").Value = -1076.2354
$ws.Range("H102").Value = 8999
$ws.Range("I102").Value = 8999
$ws.Range("K102").Value = 8999
$ws.Range("M102").Value = -7377
$ws.Range("H107").Value = 2489.3333
$ws.Range("I107").Value = 599.5
$ws.Range("J107").Value = 4001.2
$ws.Range("K107").Value = 599.5
$ws.Range("L107").Value = 4001.2
$ws.Range("M107").Value = 1320.5
$ws.Range("N107").Value = -7841.2
$ws.Range("H122").Value = 4775.107
$ws.Range("I122").Value = 3352.2
$ws.Range("J122").Value = 16632.666
$ws.Range("K122").Value = 10056.6
$ws.Range("L122").Value = 49897.99800000001
$ws.Range("M122").Value = -7606.599999999999
$ws.Range("N122").Value = -54797.99800000001
$ws.Range("H126").Value = 3115.1738
$ws.Range("I126").Value = 3041
$ws.Range("K126").Value = 9123
$ws.Range("M126").Value = -6653
$ws.Range("H132").Value = 3295181.5
$ws.Range("I132").Value = 3577183
$ws.Range("J132").Value = 5164
$ws.Range("K132").Value = 10731549
$ws.Range("L132").Value = 15492
$ws.Range("M132").Value = -10729019
$ws.Range("N132").Value = -20552

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4998
$ws.Range("I40").Value = 4998
$ws.Range("K40").Value = 4998
$ws.Range("M40").Value = -4862
$ws.Range("H43").Value = 222070
$ws.Range("I43").Value = 3745
$ws.Range("J43").Value = 309400
$ws.Range("K43").Value = 3745
$ws.Range("L43").Value = 309400
$ws.Range("M43").Value = -3552
$ws.Range("N43").Value = -309786
$ws.Range("H68").Value = 17548858
$ws.Range("I68").Value = 17548858
$ws.Range("K68").Value = 17548858
$ws.Range("M68").Value = -17548109
$ws.Range("H71").Value = 17548858
$ws.Range("I71").Value = 17548858
$ws.Range("K71").Value = 87744290
$ws.Range("M71").Value = -87740546
$ws.Range("H122").Value = 19332.334
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 14873.75
$ws.Range("I2").Value = 16665
$ws.Range("J2").Value = 9500
$ws.Range("K2").Value = 16665
$ws.Range("L2").Value = 9500
$ws.Range("M2").Value = -16553
$ws.Range("N2").Value = -9724
$ws.Range("H4").Value = 75023496
$ws.Range("I4").Value = 37500
$ws.Range("J4").Value = 150009500
$ws.Range("K4").Value = 37500
$ws.Range("L4").Value = 150009500
$ws.Range("M4").Value = -37387
$ws.Range("N4").Value = -150009726
$ws.Range("H6").Value = 100000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 100000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 100000
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = -100230
$ws.Range("H81").Value = 1449.5
$ws.Range("I81").Value = 1266.3334
$ws.Range("J81").Value = 1999
$ws.Range("K81").Value = 2532.6668
$ws.Range("L81").Value = 3998
$ws.Range("M81").Value = -1471.6668
$ws.Range("N81").Value = -6120
$ws.Range("H84").Value = 1449.5
$ws.Range("I84").Value = 1266.3334
$ws.Range("J84").Value = 1999
$ws.Range("K84").Value = 12663.334
$ws.Range("L84").Value = 19990
$ws.Range("M84").Value = -7359.333999999999
$ws.Range("N84").Value = -30598
$ws.Range("H107").Value = 996.8125
$ws.Range("I107").Value = 804.0833
$ws.Range("J107").Value = 1575
$ws.Range("K107").Value = 2412.2499
$ws.Range("L107").Value = 4725
$ws.Range("M107").Value = -492.2498999999998
$ws.Range("N107").Value = -8565
$ws.Range("H122").Value = 2971.3125
$ws.Range("I122").Value = 2972.2727
$ws.Range("J122").Value = 2969.2
$ws.Range("K122").Value = 8916.8181
$ws.Range("L122").Value = 8907.599999999999
$ws.Range("M122").Value = -6466.8181
$ws.Range("N122").Value = -13807.6
$ws.Range("H132").Value = 20005572
$ws.Range("I132").Value = 27778962
$ws.Range("K132").Value = 83336886
$ws.Range("M132").Value = -83334356
$ws.Range("H135").Value = 86353.5
$ws.Range("J135").Value = 86353.5
$ws.Range("L135").Value = 86353.5
$ws.Range("N135").Value = -96493.5
